$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2026-01 (row 26)
$ws.Range("B26").Value = 6530
$ws.Range("C26").Value = 1018
$ws.Range("D26").Value = 6082479
$ws.Range("E26").Value = 931.466921898928
$ws.Range("F26").Value = 10.35997971945242
$ws.Range("G26").Value = 8.067940552016983
$ws.Range("H26").Value = 26.66783009188214
